$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Historical Relative Freq")

# New column K: historical relative frequency = (# events occurred / # events possible) * 100
$ws.Range("K3").Formula = "=(3/29)*100"
$ws.Range("K4").Formula = "=(4/29)*100"
$ws.Range("K5").Formula = "=(5/29)*100"
$ws.Range("K6").Formula = "=(1/29)*100"
$ws.Range("K7").Formula = "=(6/29)*100"
$ws.Range("K8").Formula = "=(8/29)*100"
$ws.Range("K9").Formula = "=(2/29)*100"

$ws.Range("K3:K9").NumberFormat = "0.0"

# Footnote describing the relative-frequency calculation
$ws.Range("C12").Value = "# of event occurred / # of events it could occur"

# Leave the cursor on this sheet at G13 (matches the author's last click here)...
$ws.Range("G13").Select()

# ...then return to the workbook's originally active sheet/cell so the
# saved tab selection still points at "F24 % Cover"!C2.
$ws1 = $wb.Worksheets.Item("F24 % Cover")
$ws1.Range("C2").Select()
